$d = $word.ActiveDocument

# Update the ID placeholder text in the first paragraph, consuming the
# trailing-space run so the two runs collapse into one.
$d.Content.Find.Execute("**ID__AFFARS_pgi_5346_topic_3__ID** ", $false, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_AF_PGI_5346_202_4__ID**", 2)

# Add a paragraph border (5-twip spacing on all sides, no visible line)
# and widen the left indent on the first paragraph.
$p1 = $d.Paragraphs(1)
$p1.Range.ParagraphFormat.Borders.DistanceFromTop = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromLeft = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromBottom = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromRight = 5
$p1.Range.ParagraphFormat.LeftIndent = 11.25
